# Scheduled runner update: refresh cached market-board price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the per-job Leve
# profit sheets. Values mirror a re-pull of Universalis price data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 866.6667
$ws.Range("J32").Value = 833.3333
$ws.Range("L32").Value = 833.3333
$ws.Range("N32").Value = -1485.3333
$ws.Range("H63").Value = 28082
$ws.Range("I63").Value = 14246
$ws.Range("J63").Value = 35000
$ws.Range("K63").Value = 14246
$ws.Range("L63").Value = 35000
$ws.Range("M63").Value = -13622
$ws.Range("N63").Value = -36248
$ws.Range("H66").Value = 28082
$ws.Range("I66").Value = 14246
$ws.Range("J66").Value = 35000
$ws.Range("K66").Value = 42738
$ws.Range("L66").Value = 105000
$ws.Range("M66").Value = -39618
$ws.Range("N66").Value = -111240
$ws.Range("H88").Value = 6848.048
$ws.Range("I88").Value = 5321.3
$ws.Range("J88").Value = 8236
$ws.Range("K88").Value = 5321.3
$ws.Range("L88").Value = 8236
$ws.Range("M88").Value = -4915.3
$ws.Range("N88").Value = -9048
$ws.Range("H91").Value = 6848.048
$ws.Range("I91").Value = 5321.3
$ws.Range("J91").Value = 8236
$ws.Range("K91").Value = 5321.3
$ws.Range("L91").Value = 8236
$ws.Range("M91").Value = -3917.3
$ws.Range("N91").Value = -11044
$ws.Range("H113").Value = 4992.5
$ws.Range("I113").Value = 4571
$ws.Range("K113").Value = 4571
$ws.Range("M113").Value = -1317
$ws.Range("H116").Value = 6600.316
$ws.Range("I116").Value = 8410
$ws.Range("K116").Value = 8410
$ws.Range("M116").Value = -4968
$ws.Range("H132").Value = 5558713
$ws.Range("I132").Value = 6669599.5
$ws.Range("K132").Value = 20008798.5
$ws.Range("M132").Value = -20006268.5
$ws.Range("H137").Value = 3230037.8
$ws.Range("I137").Value = 4003971.2
$ws.Range("J137").Value = 5314.8335
$ws.Range("K137").Value = 12011913.6
$ws.Range("L137").Value = 15944.5005
$ws.Range("M137").Value = -12009363.6
$ws.Range("N137").Value = -21044.5005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7604.0796
$ws.Range("I32").Value = 5173.2104
$ws.Range("J32").Value = 22999.584
$ws.Range("K32").Value = 5173.2104
$ws.Range("L32").Value = 22999.584
$ws.Range("M32").Value = -4886.2104
$ws.Range("N32").Value = -23573.584
$ws.Range("H63").Value = 2911.7058
$ws.Range("I63").Value = 2083.25
$ws.Range("J63").Value = 4900
$ws.Range("K63").Value = 2083.25
$ws.Range("L63").Value = 4900
$ws.Range("M63").Value = -1397.25
$ws.Range("N63").Value = -6272
$ws.Range("H66").Value = 2911.7058
$ws.Range("I66").Value = 2083.25
$ws.Range("J66").Value = 4900
$ws.Range("K66").Value = 10416.25
$ws.Range("L66").Value = 24500
$ws.Range("M66").Value = -6984.25
$ws.Range("N66").Value = -31364
$ws.Range("H74").Value = 1944
$ws.Range("I74").Value = 2148.5715
$ws.Range("J74").Value = 1466.6666
$ws.Range("K74").Value = 2148.5715
$ws.Range("L74").Value = 1466.6666
$ws.Range("M74").Value = -1274.5715
$ws.Range("N74").Value = -3214.6666
$ws.Range("H77").Value = 1944
$ws.Range("I77").Value = 2148.5715
$ws.Range("J77").Value = 1466.6666
$ws.Range("K77").Value = 10742.8575
$ws.Range("L77").Value = 7333.333000000001
$ws.Range("M77").Value = -6374.8575
$ws.Range("N77").Value = -16069.333
$ws.Range("H97").Value = 537.9231
$ws.Range("I97").Value = 428.90475
$ws.Range("K97").Value = 428.90475
$ws.Range("M97").Value = 67.09525000000002
$ws.Range("H122").Value = 1594.9032
$ws.Range("I122").Value = 1043.6666
$ws.Range("J122").Value = 3484.8572
$ws.Range("K122").Value = 3130.9998
$ws.Range("L122").Value = 10454.5716
$ws.Range("M122").Value = -680.9998000000001
$ws.Range("N122").Value = -15354.5716
$ws.Range("H139").Value = 32350
$ws.Range("J139").Value = 32350
$ws.Range("L139").Value = 32350
$ws.Range("N139").Value = -42630

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 988.0833
$ws.Range("I80").Value = 861
$ws.Range("J80").Value = 1166
$ws.Range("K80").Value = 861
$ws.Range("L80").Value = 1166
$ws.Range("M80").Value = 137
$ws.Range("N80").Value = -3162
$ws.Range("H81").Value = 12560
$ws.Range("J81").Value = 12560
$ws.Range("L81").Value = 12560
$ws.Range("N81").Value = -14682
$ws.Range("H83").Value = 988.0833
$ws.Range("I83").Value = 861
$ws.Range("J83").Value = 1166
$ws.Range("K83").Value = 4305
$ws.Range("L83").Value = 5830
$ws.Range("M83").Value = 687
$ws.Range("N83").Value = -15814
$ws.Range("H84").Value = 12560
$ws.Range("J84").Value = 12560
$ws.Range("L84").Value = 37680
$ws.Range("N84").Value = -48288
$ws.Range("H86").Value = 20982.926
$ws.Range("I86").Value = 1547.95
$ws.Range("J86").Value = 76511.42999999999
$ws.Range("K86").Value = 1547.95
$ws.Range("L86").Value = 76511.42999999999
$ws.Range("M86").Value = -424.95
$ws.Range("N86").Value = -78757.42999999999
$ws.Range("H89").Value = 20982.926
$ws.Range("I89").Value = 1547.95
$ws.Range("J89").Value = 76511.42999999999
$ws.Range("K89").Value = 7739.75
$ws.Range("L89").Value = 382557.15
$ws.Range("M89").Value = -2123.75
$ws.Range("N89").Value = -393789.15
$ws.Range("H134").Value = 2957.1853
$ws.Range("I134").Value = 1867.1111
$ws.Range("J134").Value = 5137.3335
$ws.Range("K134").Value = 5601.3333
$ws.Range("L134").Value = 15412.0005
$ws.Range("M134").Value = -3066.3333
$ws.Range("N134").Value = -20482.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 12198752
$ws.Range("I58").Value = 2320.8
$ws.Range("J58").Value = 31255676
$ws.Range("K58").Value = 2320.8
$ws.Range("L58").Value = 31255676
$ws.Range("M58").Value = -2117.8
$ws.Range("N58").Value = -31256082
$ws.Range("H105").Value = 3438.75
$ws.Range("I105").Value = 2967.6667
$ws.Range("J105").Value = 4044.4285
$ws.Range("K105").Value = 2967.6667
$ws.Range("L105").Value = 4044.4285
$ws.Range("M105").Value = -1220.6667
$ws.Range("N105").Value = -7538.4285
$ws.Range("H132").Value = 2067.054
$ws.Range("I132").Value = 1832.9131
$ws.Range("J132").Value = 2451.7144
$ws.Range("K132").Value = 5498.7393
$ws.Range("L132").Value = 7355.1432
$ws.Range("M132").Value = -2968.7393
$ws.Range("N132").Value = -12415.1432
$ws.Range("H136").Value = 12198752
$ws.Range("I136").Value = 2320.8
$ws.Range("J136").Value = 31255676
$ws.Range("K136").Value = 6962.400000000001
$ws.Range("L136").Value = 93767028
$ws.Range("M136").Value = -4412.400000000001
$ws.Range("N136").Value = -93772128

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 3266.3635
$ws.Range("I123").Value = 465
$ws.Range("K123").Value = 1395
$ws.Range("M123").Value = 1055
$ws.Range("H131").Value = 2359.6296
$ws.Range("I131").Value = 1631
$ws.Range("J131").Value = 2788.2354
$ws.Range("K131").Value = 4893
$ws.Range("L131").Value = 8364.706200000001
$ws.Range("M131").Value = 147
$ws.Range("N131").Value = -18444.7062
$ws.Range("H132").Value = 1799
$ws.Range("I132").Value = 1257.8
$ws.Range("J132").Value = 2572.1428
$ws.Range("K132").Value = 11320.2
$ws.Range("L132").Value = 23149.2852
$ws.Range("M132").Value = -8790.199999999999
$ws.Range("N132").Value = -28209.2852

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3666.1082
$ws.Range("I132").Value = 3760.3333
$ws.Range("J132").Value = 3492.1538
$ws.Range("K132").Value = 11280.9999
$ws.Range("L132").Value = 10476.4614
$ws.Range("M132").Value = -8750.999899999999
$ws.Range("N132").Value = -15536.4614

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3100.111
$ws.Range("I40").Value = 2380.2
$ws.Range("K40").Value = 2380.2
$ws.Range("M40").Value = -2244.2
$ws.Range("H82").Value = 3106.6667
$ws.Range("I82").Value = 1810
$ws.Range("J82").Value = 5700
$ws.Range("K82").Value = 1810
$ws.Range("L82").Value = 5700
$ws.Range("M82").Value = -1449
$ws.Range("N82").Value = -6422
$ws.Range("H85").Value = 3106.6667
$ws.Range("I85").Value = 1810
$ws.Range("J85").Value = 5700
$ws.Range("K85").Value = 1810
$ws.Range("L85").Value = 5700
$ws.Range("M85").Value = -562
$ws.Range("N85").Value = -8196
$ws.Range("H122").Value = 2786.318
$ws.Range("I122").Value = 2339.6
$ws.Range("J122").Value = 3743.5715
$ws.Range("K122").Value = 7018.799999999999
$ws.Range("L122").Value = 11230.7145
$ws.Range("M122").Value = -4568.799999999999
$ws.Range("N122").Value = -16130.7145
$ws.Range("H132").Value = 2505.1035
$ws.Range("I132").Value = 1479.7646
$ws.Range("J132").Value = 3957.6667
$ws.Range("K132").Value = 4439.293799999999
$ws.Range("L132").Value = 11873.0001
$ws.Range("M132").Value = -1909.293799999999
$ws.Range("N132").Value = -16933.0001
$ws.Range("H136").Value = 2502233
$ws.Range("I136").Value = 3573147.5
$ws.Range("J136").Value = 3432.9167
$ws.Range("K136").Value = 10719442.5
$ws.Range("L136").Value = 10298.7501
$ws.Range("M136").Value = -10716892.5
$ws.Range("N136").Value = -15398.7501
